$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.159487842
$ws.Range("E2").Value = 1.231896994

$ws.Range("D3").Value = 1.22292089

$ws.Range("D4").Value = 1.091826347
$ws.Range("E4").Value = 1.081064855

$ws.Range("D5").Value = 1.11581984
$ws.Range("E5").Value = 1.164563028

$ws.Range("D6").Value = 1.105777922
$ws.Range("E6").Value = 1.218693595

$ws.Range("D7").Value = 1.198236411
$ws.Range("E7").Value = 1.299397493

$ws.Range("D8").Value = 1.204465398
$ws.Range("E8").Value = 1.292144173

$ws.Range("D9").Value = 1.218147487
$ws.Range("E9").Value = 1.325567315

$ws.Range("D10").Value = 1.210724827
$ws.Range("E10").Value = 1.343371159

$ws.Range("D11").Value = 1.204675149
$ws.Range("E11").Value = 1.342562654

$ws.Range("D12").Value = 1.208372733
$ws.Range("E12").Value = 1.267926445

$ws.Range("D13").Value = 1.015771636
$ws.Range("E13").Value = 1.019007469

$ws.Range("D14").Value = 1.267084355
$ws.Range("E14").Value = 1.387519745

$ws.Range("D15").Value = 1.109559743

$ws.Range("D16").Value = 1.166633613
$ws.Range("E16").Value = 1.244518276
